$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$sheet1Rows = @(
    @(6, "04:01:06", "04:02", "81_EL PELIGRO", 1),
    @(7, "04:58:02", "04:59", "11_ETCHEVERRY", 1),
    @(8, "04:44:55", "04:46", "215_EL PELIGRO", 2),
    @(9, "00:46:06", "01:12", "215_ALUAR", 26),
    @(10, "04:44:55", "05:11", "17_ROMERO", 27),
    @(11, "04:44:55", "05:21", "23_HERNANDEZ", 37),
    @(12, "04:01:06", "04:47", "215_EL PELIGRO", 46),
    @(13, "04:44:55", "05:31", "81_EL PELIGRO", 47),
    @(14, "04:44:55", "05:43", "14_ABASTO", 59),
    @(15, "03:46:12", "04:46", "215A_EL PATO", 60),
    @(16, "04:44:55", "05:51", "17_ROMERO", 67),
    @(17, "01:55:38", "03:02", "15_ABASTO", 67),
    @(18, "04:01:06", "05:12", "17_ROMERO", 71),
    @(19, "00:46:06", "01:58", "14_ABASTO", 72),
    @(20, "04:30:03", "05:44", "14_ABASTO", 74),
    @(21, "04:44:55", "06:00", "16_SANTA ANA", 76),
    @(22, "04:44:55", "06:03", "10_OLMOS", 79),
    @(23, "04:44:55", "06:10", "215A_EL PATO", 86),
    @(24, "03:46:12", "05:16", "17_ROMERO", 90),
    @(25, "04:30:03", "06:01", "16_SANTA ANA", 91),
    @(26, "04:01:06", "05:32", "81_EL PELIGRO", 91),
    @(27, "02:29:13", "04:01", "81_EL PELIGRO", 92),
    @(28, "04:58:02", "06:31", "17X38_ROMERO", 93),
    @(29, "04:58:02", "06:31", "16_SANTA ANA", 93),
    @(30, "04:30:03", "06:04", "10_OLMOS", 94),
    @(31, "01:22:42", "02:58", "215_ALUAR", 96),
    @(32, "03:46:12", "05:22", "23_HERNANDEZ", 96),
    @(33, "04:44:55", "06:23", "11_ETCHEVERRY", 99),
    @(34, "04:58:02", "06:39", "225_C ROCA-H SUR", 101),
    @(35, "04:30:03", "06:11", "215A_EL PATO", 101),
    @(36, "04:44:55", "06:26", "23_HERNANDEZ", 102),
    @(37, "04:01:06", "05:45", "14_ABASTO", 104),
    @(38, "04:44:55", "06:30", "16_SANTA ANA", 106),
    @(39, "04:44:55", "06:30", "17X38_ROMERO", 106),
    @(40, "03:46:12", "05:35", "215B_EL PATO", 109),
    @(41, "04:01:06", "05:52", "17_ROMERO", 111),
    @(42, "01:55:38", "03:48", "14_ABASTO", 113),
    @(43, "04:58:02", "06:51", "215A_EL PATO", 113),
    @(44, "03:00:53", "04:53", "11_ETCHEVERRY", 113),
    @(45, "04:44:55", "06:38", "225_C ROCA-H SUR", 114),
    @(46, "04:30:03", "06:24", "11_ETCHEVERRY", 114),
    @(47, "04:58:02", "06:54", "14_ABASTO", 116),
    @(48, "04:30:03", "06:27", "23_HERNANDEZ", 117),
    @(49, "02:47:42", "04:45", "215A_EL PATO", 118)
)

$sheet2Rows = @(
    @(6, "04:44:55", "04:46", "215_EL PELIGRO", 2),
    @(7, "00:46:06", "01:12", "215_ALUAR", 26),
    @(8, "04:01:06", "04:47", "215_EL PELIGRO", 46),
    @(9, "03:46:12", "04:46", "215A_EL PATO", 60),
    @(10, "04:44:55", "06:10", "215A_EL PATO", 86),
    @(11, "01:22:42", "02:58", "215_ALUAR", 96),
    @(12, "04:30:03", "06:11", "215A_EL PATO", 101),
    @(13, "03:46:12", "05:35", "215B_EL PATO", 109),
    @(14, "04:58:02", "06:51", "215A_EL PATO", 113),
    @(15, "02:47:42", "04:45", "215A_EL PATO", 118)
)

# Header / meta rows: Sheet1 (LP1912)
$ws1.Range("A2").Value = "Última actualización: 04:58:02"
$ws1.Range("A3").Value = "Total filas: 44"

# Header / meta rows: Sheet2 (LP1912-215)
$ws2.Range("A2").Value = "Última actualización: 04:58:02"
$ws2.Range("A3").Value = "Total filas: 10"

# Header / meta rows: Sheet3 (6203-6173)
$ws3.Range("A2").Value = "Última actualización: 04:58:02"
$ws3.Range("A3").Value = "Total filas: 44"

# Apply data rows to Sheet1 and Sheet3 (identical data tables)
foreach ($row in $sheet1Rows) {
    $r = $row[0]
    foreach ($ws in @($ws1, $ws3)) {
        $ws.Cells.Item($r, 1).Value = $row[1]
        $ws.Cells.Item($r, 2).Value = $row[2]
        $ws.Cells.Item($r, 3).Value = $row[3]
        $ws.Cells.Item($r, 4).Value = $row[4]
    }
}

# Apply data rows to Sheet2
foreach ($row in $sheet2Rows) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
}

Write-Host "Done applying schedule updates."
